$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.238.05'
$ws.Range("E2").Value = '  +3.40%  '

$ws.Range("D3").Value = '1.591.36'
$ws.Range("E3").Value = '  +1.64%  '

$ws.Range("E4").Value = '  -0.10%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '213.71'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.13%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.492'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.76%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '24.17'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +8.71%  '

$ws.Range("E9").Value = '  +0.66%  '

$ws.Range("E10").Value = '  +0.85%  '

$ws.Range("E11").Value = '  +1.71%  '

$ws.Range("D12").Value = '1.818.65'
$ws.Range("E12").Value = '  +1.66%  '

$ws.Range("D13").Value = '1.585.23'
$ws.Range("E13").Value = '  +1.37%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.532'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.23%  '

$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").Value = '28.281.06'
$ws.Range("E16").Value = '  +3.59%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '63.24'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.25%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '228.03'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +4.43%  '

$ws.Range("D19").Value = '0.0₃0710'
$ws.Range("E19").Value = '  +0.35%  '

$ws.Range("E20").Value = '  +0.34%  '

$ws.Range("E22").Value = '  -0.89%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.34'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.60%  '

$ws.Range("E24").Value = '  -0.20%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '151.96'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.34%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '15.21'
$c.Style = "Normal"

$ws.Range("E27").Value = '  +0.53%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '6.58'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.82%  '

$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  -0.35%  '

$ws.Range("E31").Value = '  +0.51%  '

$ws.Range("E32").Value = '  -0.18%  '

$ws.Range("E33").Value = '  -0.34%  '

$ws.Range("D34").Value = '1.402.98'
$ws.Range("E34").Value = '  -4.00%  '

$ws.Range("E35").Value = '  -1.84%  '

$ws.Range("E36").Value = '  -7.21%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("E38").Value = '  +0.28%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.55'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +8.59%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.543'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("E41").Value = '  -0.26%  '

$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +6.61%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.61'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -4.33%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.981'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.60%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '64.37'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.16%  '

$ws.Range("D47").Value = '1.731.24'
$ws.Range("E47").Value = '  +1.64%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '87.59'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.70%  '

$ws.Range("E49").Value = '  +0.89%  '

$ws.Range("E50").Value = '  -1.16%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0524'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
